$d = $word.ActiveDocument

# Locate the old tail range: "UbuWeb. Often, ... -->"
$old = $d.Range(315, 505)
Write-Output "OldText: [$($old.Text)]"
$tailText = $old.Text

# Collapse to end, then insert a fresh copy of the same text there
$insPoint = $d.Range(505, 505)
$insPoint.InsertAfter($tailText)

$full = $d.Content
Write-Output "Len after insert: $($full.End)"

# Delete the original (old) text range
$old2 = $d.Range(315, 505)
Write-Output "Deleting: [$($old2.Text)]"
$old2.Delete()

$full2 = $d.Content
Write-Output "Len after delete: $($full2.End)"
Write-Output $full2.Text

$ubu = $d.Range(315, 321)
Write-Output "Ubu check: [$($ubu.Text)]"
$ubu.Font.HighlightColorIndex = 3
$ubu.Font.Underline = 1
